# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# percentage-change cells for the rows whose upstream feed values moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as literal text (leading zero-padding, trailing
# zeros, and thousand-grouped dots all matter - e.g. "306.52", "0.0280").
# Prefix with an apostrophe so Excel stores the new value as text too,
# instead of silently re-typing the cell as a Number and losing formatting.
$ws.Range("D2").Value = "'42.175.09"
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = "'2.264.79"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'306.52"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = "'96.54"
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("E7").Value = '  -0.98%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").Value = "'34.99"
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = "'6.91"
$ws.Range("E13").Value = '  +2.66%  '
$ws.Range("D14").Value = "'2.616.83"
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").Value = "'14.66"
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = "'2.258.98"
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("D17").Value = "'0.791"
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = "'42.052.37"
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").Value = "'12.27"
$ws.Range("E19").Value = '  -3.30%  '
$ws.Range("D20").Value = "'0.0₃0904"
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").Value = "'67.73"
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = "'237.02"
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("D24").Value = "'2.59"
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = "'23.49"
$ws.Range("E27").Value = '  -2.75%  '
$ws.Range("D28").Value = "'37.59"
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").Value = "'162.70"
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("D32").Value = "'5.23"
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = "'3.18"
$ws.Range("E34").Value = '  +2.82%  '
$ws.Range("D35").Value = "'17.62"
$ws.Range("E35").Value = '  +2.53%  '
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("E38").Value = '  -4.32%  '
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("E41").Value = '  -3.54%  '
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("D43").Value = "'1.949.30"
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").Value = "'0.0280"
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").Value = "'54.03"
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").Value = "'92.08"
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").Value = "'71.58"
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("E51").Value = '  -2.21%  '
